$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are swapped between row 4 and row 5
$cols = @("A", "B", "D", "E", "F", "G", "H", "I")

foreach ($col in $cols) {
    $cell4 = $ws.Range($col + "4")
    $cell5 = $ws.Range($col + "5")

    $val4 = $cell4.Value2
    $val5 = $cell5.Value2

    $cell4.Value2 = $val5
    $cell5.Value2 = $val4
}
